# OLX Monitor 2026-02-22 18:30 — append newly discovered listings to the
# "PODSUMOWANIE" detail log (rows 171-178), mirroring the formatting of the
# existing last data row (170) for each column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows to append: timestamp, profile, title, price, "date added" (text,
# dd.mm.yyyy), daysOnline, url, slug
$rows = @(
    @("2026-02-22 18:30:06", "poqui", "Mieszkanie z KLIMATYZACJĄ 5 minut od UMCS, UP, KUL - Długosza", 2049, "19.12.2025", 65, "https://www.olx.pl/d/oferta/mieszkanie-z-klimatyzacja-5-minut-od-umcs-up-kul-dlugosza-CID3-ID18KAEc.html", "mieszkanie-z-klimatyzacja-5-minut-od-umcs-up-kul-dlugosza-CID3-ID18KAEc"),
    @("2026-02-22 18:30:06", "poqui", "Świeżo wykończone mieszkanie z dużym balkonem - Ponikwoda", 2299, "19.01.2026", 34, "https://www.olx.pl/d/oferta/swiezo-wykonczone-mieszkanie-z-duzym-balkonem-ponikwoda-CID3-ID1951OR.html", "swiezo-wykonczone-mieszkanie-z-duzym-balkonem-ponikwoda-CID3-ID1951OR"),
    @("2026-02-22 18:30:06", "poqui", "Kawalerka po remoncie z funkcjonalną antresolą - ul. Jana Sawy", 2499, "28.10.2025", 117, "https://www.olx.pl/d/oferta/kawalerka-po-remoncie-z-funkcjonalna-antresola-ul-jana-sawy-CID3-ID183ger.html", "kawalerka-po-remoncie-z-funkcjonalna-antresola-ul-jana-sawy-CID3-ID183ger"),
    @("2026-02-22 18:30:06", "poqui", "Przytulny pokój blisko Politechniki – ul. Przytulna", 549, "10.10.2025", 135, "https://www.olx.pl/d/oferta/przytulny-pokoj-blisko-politechniki-ul-przytulna-CID3-ID17NeTz.html", "przytulny-pokoj-blisko-politechniki-ul-przytulna-CID3-ID17NeTz"),
    @("2026-02-22 18:30:06", "pokojewlublinie", "WOLNY OD ZARAZ! Pokój jedynka, ul. Romanowskiego 58", 0, "11.08.2025", 195, "https://www.olx.pl/d/oferta/wolny-od-zaraz-pokoj-jedynka-ul-romanowskiego-58-CID3-ID16ZeYm.html", "wolny-od-zaraz-pokoj-jedynka-ul-romanowskiego-58-CID3-ID16ZeYm"),
    @("2026-02-22 18:30:06", "pokojewlublinie", "WOLNY OD ZARAZ! Super lokalizacja, blisko centrum, ul. Paganiniego 12", 12640, "19.01.2026", 34, "https://www.olx.pl/d/oferta/wolny-od-zaraz-super-lokalizacja-blisko-centrum-ul-paganiniego-12-CID3-ID195dLc.html", "wolny-od-zaraz-super-lokalizacja-blisko-centrum-ul-paganiniego-12-CID3-ID195dLc"),
    @("2026-02-22 18:30:06", "dawnypatron", "Ładny pokój jednoosobowy. Wynajmę duży pokój w centrum. ul Niecała 4.", 730, "20.09.2024", 520, "https://www.olx.pl/d/oferta/ladny-pokoj-jednoosobowy-wynajme-duzy-pokoj-w-centrum-ul-niecala-4-CID3-ID122jPM.html", "ladny-pokoj-jednoosobowy-wynajme-duzy-pokoj-w-centrum-ul-niecala-4-CID3-ID122jPM"),
    @("2026-02-22 18:30:06", "dawnypatron", "Mam do wynajęcia pokój dla os. pracującej lub studenta. Narutowicza 14", 14690, "05.12.2025", 79, "https://www.olx.pl/d/oferta/mam-do-wynajecia-pokoj-dla-os-pracujacej-lub-studenta-narutowicza-14-CID3-ID18ySfv.html", "mam-do-wynajecia-pokoj-dla-os-pracujacej-lub-studenta-narutowicza-14-CID3-ID18ySfv")
)

$lastRow = 170
$startRow = $lastRow + 1
$endRow = $startRow + $rows.Count - 1

# Replicate the formatting (styles) of the last existing data row onto the
# freshly appended rows before writing values into them.
$ws.Range("A$($lastRow):H$($lastRow)").Copy($ws.Range("A$($startRow):H$($endRow)"))

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]

    # Column E holds a "dd.mm.yyyy" string, kept as plain text in the source
    # file (not an Excel date). Some of these (day <= 12) are ambiguous and
    # Excel's smart entry would silently convert them to a date serial via
    # .Value, so force literal text with a leading apostrophe, then restore
    # the donor row's cell formatting (which clears the quote-prefix flag
    # that the apostrophe entry sets).
    $ws.Cells.Item($r, 5).Formula = "'" + $row[4]
    $ws.Cells.Item($lastRow, 5).Copy()
    $ws.Cells.Item($r, 5).PasteSpecial(-4122)

    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
}

$excel.CutCopyMode = $false
